# Applies the updated coinranking.com snapshot values to the cryptos sheet.
# Each entry is a cell plus its new display text; ForceText=$true is used for
# values that would otherwise be auto-parsed by Excel as numbers/dates, so the
# cell keeps its original plain-text storage (matching the source data feed).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "42.814.09"; ForceText = $false },
    @{ Cell = "E2"; Value = "  -6.92%  "; ForceText = $false },
    @{ Cell = "D3"; Value = "2.545.30"; ForceText = $false },
    @{ Cell = "E3"; Value = "  -1.93%  "; ForceText = $false },
    @{ Cell = "D4"; Value = "0.999"; ForceText = $true },
    @{ Cell = "E4"; Value = "  -0.02%  "; ForceText = $false },
    @{ Cell = "D5"; Value = "296.89"; ForceText = $true },
    @{ Cell = "E5"; Value = "  -4.72%  "; ForceText = $false },
    @{ Cell = "D6"; Value = "93.34"; ForceText = $true },
    @{ Cell = "E6"; Value = "  -5.18%  "; ForceText = $false },
    @{ Cell = "E7"; Value = "  -4.24%  "; ForceText = $false },
    @{ Cell = "E8"; Value = "  +0.06%  "; ForceText = $false },
    @{ Cell = "E9"; Value = "  -6.14%  "; ForceText = $false },
    @{ Cell = "D10"; Value = "35.52"; ForceText = $true },
    @{ Cell = "E10"; Value = "  -8.75%  "; ForceText = $false },
    @{ Cell = "E11"; Value = "  -3.86%  "; ForceText = $false },
    @{ Cell = "D12"; Value = "7.68"; ForceText = $true },
    @{ Cell = "E12"; Value = "  -5.42%  "; ForceText = $false },
    @{ Cell = "D13"; Value = "2.937.84"; ForceText = $false },
    @{ Cell = "E13"; Value = "  -1.74%  "; ForceText = $false },
    @{ Cell = "E14"; Value = "  -0.12%  "; ForceText = $false },
    @{ Cell = "D15"; Value = "2.545.56"; ForceText = $false },
    @{ Cell = "E15"; Value = "  -1.89%  "; ForceText = $false },
    @{ Cell = "E16"; Value = "  -5.80%  "; ForceText = $false },
    @{ Cell = "D17"; Value = "14.09"; ForceText = $true },
    @{ Cell = "E17"; Value = "  -4.83%  "; ForceText = $false },
    @{ Cell = "D18"; Value = "42.806.57"; ForceText = $false },
    @{ Cell = "E18"; Value = "  -7.05%  "; ForceText = $false },
    @{ Cell = "B19"; Value = "ShibaInu"; ForceText = $false },
    @{ Cell = "C19"; Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; ForceText = $false },
    @{ Cell = "D19"; Value = "0.0₃0975"; ForceText = $false },
    @{ Cell = "E19"; Value = "  -3.95%  "; ForceText = $false },
    @{ Cell = "B20"; Value = "Uniswap"; ForceText = $false },
    @{ Cell = "C20"; Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; ForceText = $false },
    @{ Cell = "D20"; Value = "6.62"; ForceText = $true },
    @{ Cell = "E20"; Value = "  -1.44%  "; ForceText = $false },
    @{ Cell = "D21"; Value = "12.44"; ForceText = $true },
    @{ Cell = "E21"; Value = "  -2.55%  "; ForceText = $false },
    @{ Cell = "D22"; Value = "72.11"; ForceText = $true },
    @{ Cell = "E22"; Value = "  -1.23%  "; ForceText = $false },
    @{ Cell = "D23"; Value = "259.99"; ForceText = $true },
    @{ Cell = "E23"; Value = "  -11.88%  "; ForceText = $false },
    @{ Cell = "E24"; Value = "  -4.96%  "; ForceText = $false },
    @{ Cell = "D25"; Value = "29.52"; ForceText = $true },
    @{ Cell = "E25"; Value = "  -0.04%  "; ForceText = $false },
    @{ Cell = "D26"; Value = "2.11"; ForceText = $true },
    @{ Cell = "E26"; Value = "  -5.07%  "; ForceText = $false },
    @{ Cell = "E27"; Value = "  +0.16%  "; ForceText = $false },
    @{ Cell = "D28"; Value = "10.00"; ForceText = $true },
    @{ Cell = "E28"; Value = "  -7.18%  "; ForceText = $false },
    @{ Cell = "E29"; Value = "  -4.35%  "; ForceText = $false },
    @{ Cell = "D30"; Value = "35.94"; ForceText = $true },
    @{ Cell = "E30"; Value = "  -6.56%  "; ForceText = $false },
    @{ Cell = "D31"; Value = "5.94"; ForceText = $true },
    @{ Cell = "E31"; Value = "  -4.36%  "; ForceText = $false },
    @{ Cell = "D32"; Value = "150.87"; ForceText = $true },
    @{ Cell = "E32"; Value = "  -3.01%  "; ForceText = $false },
    @{ Cell = "E33"; Value = "  -1.99%  "; ForceText = $false },
    @{ Cell = "E34"; Value = "  -5.54%  "; ForceText = $false },
    @{ Cell = "E35"; Value = "  -2.76%  "; ForceText = $false },
    @{ Cell = "E36"; Value = "  -5.34%  "; ForceText = $false },
    @{ Cell = "E37"; Value = "  -6.75%  "; ForceText = $false },
    @{ Cell = "E38"; Value = "  +14.28%  "; ForceText = $false },
    @{ Cell = "E39"; Value = "  -3.42%  "; ForceText = $false },
    @{ Cell = "D40"; Value = "16.10"; ForceText = $true },
    @{ Cell = "E40"; Value = "  +2.59%  "; ForceText = $false },
    @{ Cell = "D41"; Value = "3.41"; ForceText = $true },
    @{ Cell = "E41"; Value = "  -4.34%  "; ForceText = $false },
    @{ Cell = "E42"; Value = "  -6.21%  "; ForceText = $false },
    @{ Cell = "D43"; Value = "2.079.49"; ForceText = $false },
    @{ Cell = "E43"; Value = "  -1.24%  "; ForceText = $false },
    @{ Cell = "E44"; Value = "  -3.45%  "; ForceText = $false },
    @{ Cell = "E45"; Value = "  +0.04%  "; ForceText = $false },
    @{ Cell = "D46"; Value = "85.00"; ForceText = $true },
    @{ Cell = "E46"; Value = "  -13.60%  "; ForceText = $false },
    @{ Cell = "B47"; Value = "ApeXProtocol"; ForceText = $false },
    @{ Cell = "C47"; Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"; ForceText = $false },
    @{ Cell = "D47"; Value = "1.58"; ForceText = $true },
    @{ Cell = "E47"; Value = "  +2.95%  "; ForceText = $false },
    @{ Cell = "B48"; Value = "RocketPoolETH"; ForceText = $false },
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"; ForceText = $false },
    @{ Cell = "D48"; Value = "2.791.38"; ForceText = $false },
    @{ Cell = "E48"; Value = "  -1.94%  "; ForceText = $false },
    @{ Cell = "B49"; Value = "Stacks"; ForceText = $false },
    @{ Cell = "C49"; Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"; ForceText = $false },
    @{ Cell = "D49"; Value = "1.70"; ForceText = $true },
    @{ Cell = "E49"; Value = "  -1.70%  "; ForceText = $false },
    @{ Cell = "D50"; Value = "103.87"; ForceText = $true },
    @{ Cell = "E50"; Value = "  -4.17%  "; ForceText = $false },
    @{ Cell = "D51"; Value = "8.66"; ForceText = $true },
    @{ Cell = "E51"; Value = "  -9.61%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Force text storage, write the value, then restore the default
        # (unstyled) cell style so only the content changes.
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.NumberFormat = "@"
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
